$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ciudades")

# Update the "last updated" timestamp text (A1)
$ws.Range("A1").Value = "Datos actualizados a 25 de Marzo de 2020 a las 22:46"

# Update Cataluña row (row 5) figures
$ws.Range("B5").Value = 11592
$ws.Range("C5").Value = 1697
$ws.Range("D5").Value = 9223
$ws.Range("E5").Value = 672
